$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.244.55'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').Value = '2.579.96'
$ws.Range('E3').Value = '  -1.11%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.72%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.49%  '

$ws.Range('E7').Value = '  +0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.595'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.29%  '

$ws.Range('D9').Value = '2.584.27'
$ws.Range('E9').Value = '  -1.54%  '

$ws.Range('E10').Value = '  -1.45%  '

$ws.Range('E11').Value = '  +2.90%  '

$ws.Range('E12').Value = '  +11.00%  '

$ws.Range('E13').Value = '  +2.68%  '

$ws.Range('D14').Value = '3.031.48'
$ws.Range('E14').Value = '  -1.15%  '

$ws.Range('D15').Value = '59.249.83'
$ws.Range('E15').Value = '  +0.50%  '

$ws.Range('E16').Value = '  +6.49%  '

$ws.Range('E17').Value = '  +3.55%  '

$ws.Range('D18').Value = '2.588.18'
$ws.Range('E18').Value = '  -1.66%  '

$ws.Range('E19').Value = '  +1.54%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '338.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.39%  '

$ws.Range('E21').Value = '  +1.41%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.85%  '

$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.28%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.462'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.78%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('E27').Value = '  -0.36%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.29'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.68%  '

$ws.Range('D29').Value = '0.0₃0784'
$ws.Range('E29').Value = '  +2.59%  '

$ws.Range('E31').Value = '  +0.51%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.94%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.55%  '

$ws.Range('E34').Value = '  +0.44%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.61%  '

$ws.Range('E36').Value = '  +2.32%  '

$ws.Range('E37').Value = '  -3.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.872'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.32%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.48%  '

$ws.Range('E40').Value = '  +1.43%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '296.52'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.40%  '

$ws.Range('E42').Value = '  +1.74%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.10%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '131.23'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +11.54%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0978'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.22%  '

$ws.Range('E46').Value = '  -1.27%  '

$ws.Range('E47').Value = '  -0.50%  '

$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.66'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.17%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.22'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.93%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0234'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.86%  '

$ws.Range('D51').Value = '1.957.00'
$ws.Range('E51').Value = '  +0.38%  '
